$p = $ppt.ActivePresentation

# --- 1. Update the "update automatically" date placeholder cached text
#        (slide master + every custom layout) from 5/9/13 to 7/4/14 ---
function Set-DatePlaceholderText($shapes, $text) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.HasTextFrame -and $shp.PlaceholderFormat.Type -eq 16) {
            $shp.TextFrame.TextRange.Text = $text
        }
    }
}

$master = $p.SlideMaster
Set-DatePlaceholderText $master.Shapes "7/4/14"
for ($j = 1; $j -le $master.CustomLayouts.Count; $j++) {
    $layout = $master.CustomLayouts.Item($j)
    Set-DatePlaceholderText $layout.Shapes "7/4/14"
}

# --- 2. Fix Randy's byline: "Randall J. Pruim" -> "Randall Pruim" ---
$s = $p.Slides.Item(1)
$nameShape = $s.Shapes.Item(1)
$nameRange = $nameShape.TextFrame.TextRange
$byline = $nameRange.Paragraphs(3, 1)
$byline.Text = "Randall Pruim"
